$wb = $excel.ActiveWorkbook

# --- Step 1: Create the new "2022-Q4" sheet by copying the "2022-Q2" sheet
#     (keeps the same column layout/styles), positioned right before "2022-Q2".
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ4 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ4.Name = "2022-Q4"

# --- Step 2: Clear the old data rows (2-10) from the copy, keep header row 1.
$wsQ4.Range("A2:A10").EntireRow.Delete()

# Restore the index-column style (col A uses the same style as the header
# cells) on the two rows we are about to (re)populate.
$wsQ4.Cells.Item(1, 2).Copy()
$wsQ4.Cells.Item(2, 1).PasteSpecial(-4122)
$wsQ4.Cells.Item(3, 1).PasteSpecial(-4122)

# Columns B-G hold text values (fund codes with leading zeros, decimal
# strings, etc.) - force text format before assigning so they are not
# auto-coerced to numbers.
$wsQ4.Range("B2:G3").NumberFormat = "@"

# --- Step 3: Fill the new 2022-Q4 fund holdings data.
$wsQ4.Cells.Item(2, 1).Value = 0
$wsQ4.Cells.Item(2, 2).Value = "014668"
$wsQ4.Cells.Item(2, 3).Value = "银华专精特新量化优选股票A"
$wsQ4.Cells.Item(2, 4).Value = "0.26"
$wsQ4.Cells.Item(2, 5).Value = "94.19"
$wsQ4.Cells.Item(2, 6).Value = "1.31"
$wsQ4.Cells.Item(2, 7).Value = "0.0034"
$wsQ4.Cells.Item(2, 8).Value = 10

$wsQ4.Cells.Item(3, 1).Value = 1
$wsQ4.Cells.Item(3, 2).Value = "014669"
$wsQ4.Cells.Item(3, 3).Value = "银华专精特新量化优选股票C"
$wsQ4.Cells.Item(3, 4).Value = "0.15"
$wsQ4.Cells.Item(3, 5).Value = "94.19"
$wsQ4.Cells.Item(3, 6).Value = "1.31"
$wsQ4.Cells.Item(3, 7).Value = "0.0020"
$wsQ4.Cells.Item(3, 8).Value = 10

# Drop the temporary text-format styling so these cells end up with the
# default (unstyled) cell format, matching the rest of the data rows.
$wsQ4.Range("B2:G3").ClearFormats()

# --- Step 4: Update the "总计" (summary) sheet - insert a new row for
#     2022-Q4 above the existing 2022-Q2 row, shifting the rest down.
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(2, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.01
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2

Write-Host "Edit complete"
